$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Copy()
$ws.Range("D3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D3").Value = "IN_STOCK"

$ws.Range("E4").Value = "Сотрудник 3"

$ws.Range("E3:E4").Select()
